# CricDream issue list -- add "9Sep2020" sheet, hide old Sheet2, update Sheet1 view.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet1: move selection/frozen pane before we switch the active sheet away
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A2:F4").Select()

# ---------------------------------------------------------------------------
# 2. Insert the new "9Sep2020" worksheet right after Sheet1
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "9Sep2020"

# ---------------------------------------------------------------------------
# 3. Populate the new worksheet
# ---------------------------------------------------------------------------

# Title row
$ws2.Range("A1").Value = "Admin login test (apurva) on 9th September 2020"
$ws2.Range("A1:C1").Merge()
$ws2.Range("A1").Font.Bold = $true
$ws2.Range("A1").Font.Size = 16

# Header row
$ws2.Range("A2").Value = "Sr. No."
$ws2.Range("B2").Value = "Issue Description"
$ws2.Range("C2").Value = "Priority"
$ws2.Range("D2").Value = "Issue Date"
$ws2.Range("E2").Value = "Status"
$ws2.Range("F2").Value = "Resolved Date"
$ws2.Range("G2").Value = "Resolved Description"
$ws2.Range("A2:G2").Font.Bold = $true

# Issue rows
$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = "Sold button to be disabled, if already purchased"
$ws2.Range("C3").Value = "low"
$ws2.Range("D3").Value = "9/9/2020"
$ws2.Range("E3").Value = "Pending"

$ws2.Range("A4").Value = 2
$ws2.Range("B4").Value = "Reconfirm message for sold"
$ws2.Range("C4").Value = "High"
$ws2.Range("D4").Value = "9/7/2020"
$ws2.Range("E4").Value = "Pending"

$ws2.Range("A5").Value = 3
$ws2.Range("B5").Value = "RHS to automatic hide after confirmed?"
$ws2.Range("C5").Value = "High"

$ws2.Range("A6").Value = 4
$ws2.Range("B6").Value = "Cancel button in RHS along with Confirm"
$ws2.Range("C6").Value = "low"

$ws2.Range("A7").Value = 5
$ws2.Range("B7").Value = "Stat to set to IPL2020"
$ws2.Range("C7").Value = "high"

$ws2.Range("A8").Value = 6
$ws2.Range("B8").Value = "Admin to get Auction start button"
$ws2.Range("C8").Value = "low"
$ws2.Range("I8").Value = "Done"

$ws2.Range("A9").Value = 7
$ws2.Range("B9").Value = "Unsold to be implemented"
$ws2.Range("C9").Value = "high"

$ws2.Range("A10").Value = 8
$ws2.Range("B10").Value = "Logout to be implemented"
$ws2.Range("C10").Value = "low"

$ws2.Range("A11").Value = 9
$ws2.Range("B11").Value = "Check if (as user) My team shows purchased players"
$ws2.Range("C11").Value = "high"
$ws2.Range("I11").Value = "Done"

$ws2.Range("A12").Value = 10
$ws2.Range("B12").Value = "(admin) My team do not get all players sold"
$ws2.Range("C12").Value = "high"

$ws2.Range("A13").Value = 11
$ws2.Range("B13").Value = "Extra Menu options to be added"

$ws2.Range("A14").Value = 12
$ws2.Range("B14").Value = "Balance not getting updated for franchisee"

$ws2.Range("A15").Value = 13
$ws2.Range("B15").Value = "Auction sequnce from PENDING to RUNNING working fine"

# ---------------------------------------------------------------------------
# 4. Formatting: borders around the used area, centered Sr.No/Priority, widths
# ---------------------------------------------------------------------------
$ws2.Range("A2:G15").Borders.LineStyle = 1
$ws2.Range("A2:A15").HorizontalAlignment = -4108
$ws2.Range("C2:C15").HorizontalAlignment = -4108
$ws2.Range("B3:B15").WrapText = $true

$ws2.Columns.Item(1).ColumnWidth = 8.85546875
$ws2.Columns.Item(2).ColumnWidth = 53.28515625
$ws2.Columns.Item(3).ColumnWidth = 9.85546875
$ws2.Columns.Item(4).ColumnWidth = 20.28515625
$ws2.Columns.Item(5).ColumnWidth = 13.28515625
$ws2.Columns.Item(6).ColumnWidth = 18.5703125
$ws2.Columns.Item(8).ColumnWidth = 20.140625

$ws2.Columns.Item(4).Hidden = $true
$ws2.Columns.Item(5).Hidden = $true
$ws2.Columns.Item(6).Hidden = $true
$ws2.Columns.Item(7).Hidden = $true
$ws2.Columns.Item(8).Hidden = $true

$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Data validation dropdown for Status, sourced from the (now hidden) Sheet2
$ws2.Range("E3:E4").Validation.Add(3, 1, 1, "=Sheet2!$A$1:$A$10")

# ---------------------------------------------------------------------------
# 5. Hide the old Sheet2 (Pending/Resolved lookup list) -- content untouched
# ---------------------------------------------------------------------------
$wsOld = $wb.Worksheets.Item("Sheet2")
$wsOld.Visible = $false

# ---------------------------------------------------------------------------
# 6. Make the new sheet the active tab
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("K11").Select()
